$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "PO"
$ws.Range("B1").Value = "Date"
$ws.Range("C1").Value = "CurrencyRate"
$ws.Range("D1").Value = "Amount"

# Data row
$ws.Range("A2").Value = "testest"
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Value = (Get-Date -Year 2025 -Month 3 -Day 21 -Hour 0 -Minute 0 -Second 0).Date
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 0

$ws.Range("E6").Select()
